$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# HU - Modulo Usuario Profesor progress updated from 30% to 80%
$ws.Range("H22").Value = 0.8

# Clear the "Sprint retrospective" progress cells that no longer have a value
$ws.Range("H19").ClearContents()
$ws.Range("H23").ClearContents()

# Update the Hito 1 milestone average formula (now averages 8 items instead of 9,
# since H19 no longer contributes a numeric value)
$ws.Range("H20").Formula = "=(H10+H11+H12+H14+H15+H16+H17+H18)/8"

# Update the Hito 2 milestone average formula (now averages 2 items instead of 3,
# since H23 no longer contributes a numeric value)
$ws.Range("H24").Formula = "=(H21+H22)/2"

# Leave the selection where the author left it when saving
$ws.Range("I23").Select()

